# Auto-generated script: updates Leve profit-calculation values across all sheets
# to match the scheduled-runner price refresh recorded in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1880.5
$ws.Range("J125").Value = 1496
$ws.Range("L125").Value = 13464
$ws.Range("N125").Value = -18384
$ws.Range("H132").Value = 24334.904
$ws.Range("I132").Value = 3120.8057
$ws.Range("J132").Value = 151619.5
$ws.Range("K132").Value = 9362.417099999999
$ws.Range("L132").Value = 454858.5
$ws.Range("M132").Value = -6832.417099999999
$ws.Range("N132").Value = -459918.5
$ws.Range("H137").Value = 6752.65
$ws.Range("I137").Value = 8127.294
$ws.Range("J137").Value = 5736.609
$ws.Range("K137").Value = 24381.882
$ws.Range("L137").Value = 17209.827
$ws.Range("M137").Value = -21831.882
$ws.Range("N137").Value = -22309.827

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 33334662
$ws.Range("I45").Value = 41667670
$ws.Range("J45").Value = 2619
$ws.Range("K45").Value = 41667670
$ws.Range("L45").Value = 2619
$ws.Range("M45").Value = -41667293
$ws.Range("N45").Value = -3373
$ws.Range("H61").Value = 2569.9666
$ws.Range("I61").Value = 1990.9412
$ws.Range("K61").Value = 1990.9412
$ws.Range("M61").Value = -1778.9412
$ws.Range("H132").Value = 19234354
$ws.Range("I132").Value = 45457236
$ws.Range("J132").Value = 4241.8
$ws.Range("K132").Value = 136371708
$ws.Range("L132").Value = 12725.4
$ws.Range("M132").Value = -136369178
$ws.Range("N132").Value = -17785.4
$ws.Range("H136").Value = 2569.9666
$ws.Range("I136").Value = 1990.9412
$ws.Range("K136").Value = 5972.8236
$ws.Range("M136").Value = -3422.8236

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2870.4614
$ws.Range("I134").Value = 2391
$ws.Range("J134").Value = 3949.25
$ws.Range("K134").Value = 7173
$ws.Range("L134").Value = 11847.75
$ws.Range("M134").Value = -4638
$ws.Range("N134").Value = -16917.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4834740
$ws.Range("I31").Value = 1642.28
$ws.Range("J31").Value = 7580818
$ws.Range("K31").Value = 1642.28
$ws.Range("L31").Value = 7580818
$ws.Range("M31").Value = -1347.28
$ws.Range("N31").Value = -7581408
$ws.Range("H34").Value = 4834740
$ws.Range("I34").Value = 1642.28
$ws.Range("J34").Value = 7580818
$ws.Range("K34").Value = 1642.28
$ws.Range("L34").Value = 7580818
$ws.Range("M34").Value = -1440.28
$ws.Range("N34").Value = -7581222
$ws.Range("H58").Value = 12501523
$ws.Range("I58").Value = 892.46155
$ws.Range("J58").Value = 35716980
$ws.Range("K58").Value = 892.46155
$ws.Range("L58").Value = 35716980
$ws.Range("M58").Value = -689.46155
$ws.Range("N58").Value = -35717386
$ws.Range("H132").Value = 56764.883
$ws.Range("I132").Value = 2404.5881
$ws.Range("J132").Value = 159445.44
$ws.Range("K132").Value = 7213.7643
$ws.Range("L132").Value = 478336.32
$ws.Range("M132").Value = -4683.7643
$ws.Range("N132").Value = -483396.32
$ws.Range("H134").Value = 342618.16
$ws.Range("I134").Value = 979.0833
$ws.Range("J134").Value = 2802419.5
$ws.Range("K134").Value = 2937.2499
$ws.Range("L134").Value = 8407258.5
$ws.Range("M134").Value = -402.2498999999998
$ws.Range("N134").Value = -8412328.5
$ws.Range("H136").Value = 12501523
$ws.Range("I136").Value = 892.46155
$ws.Range("J136").Value = 35716980
$ws.Range("K136").Value = 2677.38465
$ws.Range("L136").Value = 107150940
$ws.Range("M136").Value = -127.38465
$ws.Range("N136").Value = -107156040

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1966.6666
$ws.Range("J75").Value = 2450
$ws.Range("L75").Value = 7350
$ws.Range("N75").Value = -9346
$ws.Range("H78").Value = 1966.6666
$ws.Range("J78").Value = 2450
$ws.Range("L78").Value = 22050
$ws.Range("N78").Value = -32034
$ws.Range("H92").Value = 2250
$ws.Range("I92").Value = 2000
$ws.Range("J92").Value = 3000
$ws.Range("K92").Value = 6000
$ws.Range("L92").Value = 9000
$ws.Range("M92").Value = -4752
$ws.Range("N92").Value = -11496
$ws.Range("H94").Value = 500
$ws.Range("I94").Value = 500
$ws.Range("K94").Value = 1500
$ws.Range("M94").Value = -824
$ws.Range("H95").Value = 1335000
$ws.Range("J95").Value = 1335000
$ws.Range("L95").Value = 4005000
$ws.Range("N95").Value = -4009118
$ws.Range("H96").Value = 4205.6
$ws.Range("J96").Value = 4205.6
$ws.Range("L96").Value = 12616.8
$ws.Range("N96").Value = -16734.8
$ws.Range("H97").Value = 81760
$ws.Range("I97").Value = 300
$ws.Range("J97").Value = 102125
$ws.Range("K97").Value = 900
$ws.Range("L97").Value = 306375
$ws.Range("M97").Value = -404
$ws.Range("N97").Value = -307367
$ws.Range("H98").Value = 141
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("H99").Value = 2899.8
$ws.Range("I99").Value = 1633
$ws.Range("K99").Value = 4899
$ws.Range("M99").Value = -2653
$ws.Range("H100").Value = 4983.3335
$ws.Range("J100").Value = 5580
$ws.Range("L100").Value = 16740
$ws.Range("N100").Value = -18362
$ws.Range("H101").Value = 8023.2
$ws.Range("J101").Value = 8023.2
$ws.Range("L101").Value = 24069.6
$ws.Range("N101").Value = -28937.6
$ws.Range("H102").Value = 8000
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("H103").Value = 2733.3333
$ws.Range("I103").Value = 1600
$ws.Range("J103").Value = 5000
$ws.Range("K103").Value = 4800
$ws.Range("L103").Value = 15000
$ws.Range("M103").Value = -3921
$ws.Range("N103").Value = -16758
$ws.Range("H105").Value = 115577.336
$ws.Range("J105").Value = 115577.336
$ws.Range("L105").Value = 346732.008
$ws.Range("N105").Value = -351974.008
$ws.Range("H106").Value = 869666.7
$ws.Range("I106").Value = 2600000
$ws.Range("J106").Value = 4500
$ws.Range("K106").Value = 7800000
$ws.Range("L106").Value = 13500
$ws.Range("M106").Value = -7799054
$ws.Range("N106").Value = -15392
$ws.Range("H109").Value = 2213.5107
$ws.Range("J109").Value = 3435.3845
$ws.Range("L109").Value = 10306.1535
$ws.Range("N109").Value = -12386.1535
$ws.Range("H115").Value = 4199.8
$ws.Range("I115").Value = 428
$ws.Range("J115").Value = 4780.077
$ws.Range("K115").Value = 1284
$ws.Range("L115").Value = 14340.231
$ws.Range("M115").Value = -109
$ws.Range("N115").Value = -16690.231
$ws.Range("H118").Value = 3652.6924
$ws.Range("I118").Value = 543
$ws.Range("J118").Value = 4585.6
$ws.Range("K118").Value = 1629
$ws.Range("L118").Value = 13756.8
$ws.Range("M118").Value = -386
$ws.Range("N118").Value = -16242.8
$ws.Range("H132").Value = 3293.625
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 3549.8572
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 31948.7148
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -37008.7148
$ws.Range("N98").ClearContents()
$ws.Range("M102").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 98975
$ws.Range("J62").Value = 98975
$ws.Range("L62").Value = 98975
$ws.Range("N62").Value = -100347
$ws.Range("H65").Value = 98975
$ws.Range("J65").Value = 98975
$ws.Range("L65").Value = 296925
$ws.Range("N65").Value = -303789
$ws.Range("H132").Value = 3614.35
$ws.Range("I132").Value = 2639.9
$ws.Range("J132").Value = 4588.8
$ws.Range("K132").Value = 7919.700000000001
$ws.Range("L132").Value = 13766.4
$ws.Range("M132").Value = -5389.700000000001
$ws.Range("N132").Value = -18826.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 66670570
$ws.Range("I7").Value = 111113720
$ws.Range("J7").Value = 5839.6665
$ws.Range("K7").Value = 111113720
$ws.Range("L7").Value = 5839.6665
$ws.Range("M7").Value = -111113608
$ws.Range("N7").Value = -6063.6665
$ws.Range("H16").Value = 3387.1
$ws.Range("I16").Value = 3271.5
$ws.Range("J16").Value = 3849.5
$ws.Range("K16").Value = 3271.5
$ws.Range("L16").Value = 3849.5
$ws.Range("M16").Value = -3101.5
$ws.Range("N16").Value = -4189.5
$ws.Range("H40").Value = 6834.1113
$ws.Range("I40").Value = 7101
$ws.Range("J40").Value = 6500.5
$ws.Range("K40").Value = 7101
$ws.Range("L40").Value = 6500.5
$ws.Range("M40").Value = -6965
$ws.Range("N40").Value = -6772.5
$ws.Range("H126").Value = 66670570
$ws.Range("I126").Value = 111113720
$ws.Range("J126").Value = 5839.6665
$ws.Range("K126").Value = 333341160
$ws.Range("L126").Value = 17518.9995
$ws.Range("M126").Value = -333338690
$ws.Range("N126").Value = -22458.9995
$ws.Range("H132").Value = 4493.75
$ws.Range("I132").Value = 3510.2222
$ws.Range("J132").Value = 5298.4546
$ws.Range("K132").Value = 10530.6666
$ws.Range("L132").Value = 15895.3638
$ws.Range("M132").Value = -8000.6666
$ws.Range("N132").Value = -20955.3638
$ws.Range("H136").Value = 2018.4286
$ws.Range("I136").Value = 1491.35
$ws.Range("J136").Value = 3336.125
$ws.Range("K136").Value = 4474.049999999999
$ws.Range("L136").Value = 10008.375
$ws.Range("M136").Value = -1924.049999999999
$ws.Range("N136").Value = -15108.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2298.84
$ws.Range("I132").Value = 1505.3572
$ws.Range("J132").Value = 3308.7273
$ws.Range("K132").Value = 4516.071599999999
$ws.Range("L132").Value = 9926.1819
$ws.Range("M132").Value = -1986.071599999999
$ws.Range("N132").Value = -14986.1819

